$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @{ A="ECs"; B="Vcam1"; C="Itgb2"; D="ECs"; E=3; F=1; G=20.408218; H=61.224654; I=0.1108535210972707; J=0.1108535210972707; K=1; L=0.3333333333333333; M=0.1145113333333333; N=0.343534; O=0.001785365609625045; P=0.001785365609625044; Q=2.336972254137334; R=21.032750287236; S=0.0001979140642729114; T=0.0001979140642729114 }
  @{ A="ECs"; B="Vcam1"; C="Itgb2"; D="FAPs"; E=3; F=1; G=20.408218; H=61.224654; I=0.1108535210972707; J=0.1108535210972707; K=3; L=1; M=0.467525; N=1.402575; O=0.007289261528465441; P=0.007289261528465441; Q=9.541352120450002; R=85.87216908405001; S=0.0008080403066292674; T=0.0008080403066292674 }
  @{ A="ECs"; B="Vcam1"; C="Itgb2"; D="Resolving-Mac"; E=3; F=1; G=20.408218; H=61.224654; I=0.1108535210972707; J=0.1108535210972707; K=3; L=1; M=63.556834; N=190.670502; O=0.9909253728619096; P=0.9909253728619095; Q=1297.081723661812; R=11673.73551295631; S=0.1098475667263685; T=0.1098475667263685 }
  @{ A="FAPs"; B="Vcam1"; C="Itgb2"; D="ECs"; E=3; F=1; G=47.25592399999999; H=141.767772; I=0.2566851044076959; J=0.256685104407696; K=1; L=0.3333333333333333; M=0.1145113333333333; N=0.343534; O=0.001785365609625045; P=0.001785365609625044; Q=5.411338865138666; R=48.702049786248; S=0.0004582767579125142; T=0.0004582767579125143 }
  @{ A="FAPs"; B="Vcam1"; C="Itgb2"; D="FAPs"; E=3; F=1; G=47.25592399999999; H=141.767772; I=0.2566851044076959; J=0.256685104407696; K=3; L=1; M=0.467525; N=1.402575; O=0.007289261528465441; P=0.007289261528465441; Q=22.0933258681; R=198.8399328129; S=0.001871044856489153; T=0.001871044856489153 }
  @{ A="FAPs"; B="Vcam1"; C="Itgb2"; D="Resolving-Mac"; E=3; F=1; G=47.25592399999999; H=141.767772; I=0.2566851044076959; J=0.256685104407696; K=3; L=1; M=63.556834; N=190.670502; O=0.9909253728619096; P=0.9909253728619095; Q=3003.436917184616; R=27030.93225466154; S=0.2543557827932943; T=0.2543557827932943 }
  @{ A="MuSCs"; B="Vcam1"; C="Itgb2"; D="ECs"; E=3; F=1; G=85.307233; H=255.921699; I=0.4633725077375833; J=0.4633725077375833; K=1; L=0.3333333333333333; M=0.1145113333333333; N=0.343534; O=0.001785365609625045; P=0.001785365609625044; Q=9.768644993807333; R=87.917804944266; S=0.0008272893397603961; T=0.0008272893397603961 }
  @{ A="MuSCs"; B="Vcam1"; C="Itgb2"; D="FAPs"; E=3; F=1; G=85.307233; H=255.921699; I=0.4633725077375833; J=0.4633725077375833; K=3; L=1; M=0.467525; N=1.402575; O=0.007289261528465441; P=0.007289261528465441; Q=39.883264108325; R=358.949376974925; S=0.00337764339400012; T=0.003377643394000121 }
  @{ A="MuSCs"; B="Vcam1"; C="Itgb2"; D="Resolving-Mac"; E=3; F=1; G=85.307233; H=255.921699; I=0.4633725077375833; J=0.4633725077375833; K=3; L=1; M=63.556834; N=190.670502; O=0.9909253728619096; P=0.9909253728619095; Q=5421.857646780322; R=48796.71882102289; S=0.4591675750038228; T=0.4591675750038228 }
  @{ A="Resolving-Mac"; B="Vcam1"; C="Itgb2"; D="ECs"; E=3; F=1; G=31.12938966666666; H=93.38816899999999; I=0.16908886675745; J=0.16908886675745; K=1; L=0.3333333333333333; M=0.1145113333333333; N=0.343534; O=0.001785365609625045; P=0.001785365609625044; Q=3.564667916582889; R=32.08201124924599; S=0.0003018854476792227; T=0.0003018854476792227 }
  @{ A="Resolving-Mac"; B="Vcam1"; C="Itgb2"; D="FAPs"; E=3; F=1; G=31.12938966666666; H=93.38816899999999; I=0.16908886675745; J=0.16908886675745; K=3; L=1; M=0.467525; N=1.402575; O=0.007289261528465441; P=0.007289261528465441; Q=14.55376790390833; R=130.983911135175; S=0.001232532971346899; T=0.001232532971346899 }
  @{ A="Resolving-Mac"; B="Vcam1"; C="Itgb2"; D="Resolving-Mac"; E=3; F=1; G=31.12938966666666; H=93.38816899999999; I=0.16908886675745; J=0.16908886675745; K=3; L=1; M=63.556834; N=190.670502; O=0.9909253728619096; P=0.9909253728619095; Q=1978.485451565649; R=17806.36906409084; S=0.1675544483384239; T=0.1675544483384239 }
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r++
}

# Remove the trailing 4 rows (14-17) that no longer apply (MuSCs target removed)
$ws.Range("A14:T17").EntireRow.Delete() | Out-Null

Write-Output "done"